$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=2; D='67.678.79'; E='  -2.57%  '},
    @{Row=3; D='3.798.36'; E='  +0.90%  '},
    @{Row=4; E='  -0.01%  '},
    @{Row=5; D='601.38'; E='  -2.25%  '; DText=$true},
    @{Row=6; D='171.99'; E='  -3.17%  '; DText=$true},
    @{Row=7; D='3.798.19'; E='  +0.92%  '},
    @{Row=8; E='  +0.05%  '},
    @{Row=9; D='0.531'; E='  +0.61%  '; DText=$true},
    @{Row=10; E='  -4.94%  '},
    @{Row=11; D='6.24'; E='  -6.50%  '; DText=$true},
    @{Row=12; D='0.465'; E='  -3.73%  '; DText=$true},
    @{Row=13; D='38.52'; E='  -3.73%  '; DText=$true},
    @{Row=14; D='0.0000243'; E='  -3.72%  '; DText=$true},
    @{Row=15; D='4.436.87'; E='  +0.91%  '},
    @{Row=16; D='3.804.64'; E='  +0.92%  '},
    @{Row=17; D='67.707.48'; E='  -2.59%  '},
    @{Row=18; D='7.21'; E='  -4.02%  '; DText=$true},
    @{Row=19; E='  -3.94%  '},
    @{Row=20; E='  +3.86%  '},
    @{Row=21; D='491.76'; E='  -3.34%  '; DText=$true},
    @{Row=22; D='9.18'; E='  -2.06%  '; DText=$true},
    @{Row=23; D='0.741'; E='  +1.93%  '; DText=$true},
    @{Row=24; D='85.69'; E='  -0.80%  '; DText=$true},
    @{Row=25; B='Fetch.AI'; C='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D='2.37'; E='  -5.49%  '; DText=$true},
    @{Row=26; B='PEPE'; C='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D='0.0000145'; E='  +7.19%  '; DText=$true},
    @{Row=27; D='12.29'; E='  -3.93%  '; DText=$true},
    @{Row=28; D='10.22'; E='  -3.11%  '; DText=$true},
    @{Row=29; E='  +0.00%  '},
    @{Row=30; D='2.97'; E='  +0.69%  '; DText=$true},
    @{Row=31; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='2.43'; E='  -2.73%  '; DText=$true},
    @{Row=32; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='32.82'; E='  +7.01%  '; DText=$true},
    @{Row=33; D='7.78'; E='  -2.45%  '; DText=$true},
    @{Row=34; D='0.109'; E='  -4.13%  '; DText=$true},
    @{Row=35; E='  +0.00%  '},
    @{Row=36; D='1.01'; E='  -3.84%  '; DText=$true},
    @{Row=37; D='5.80'; E='  -5.30%  '; DText=$true},
    @{Row=38; B='Kaspa'; C='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D='0.132'; E='  -4.97%  '; DText=$true},
    @{Row=39; B='TheGraph'; C='https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; D='0.328'; E='  -3.42%  '; DText=$true},
    @{Row=40; B='Bittensor'; C='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D='456.39'; E='  +1.32%  '; DText=$true},
    @{Row=41; D='49.03'; E='  -1.71%  '; DText=$true},
    @{Row=42; D='2.01'; E='  -2.94%  '; DText=$true},
    @{Row=43; D='2.85'; E='  -4.19%  '; DText=$true},
    @{Row=44; D='8.42'; E='  -1.44%  '; DText=$true},
    @{Row=45; D='41.16'; E='  -7.82%  '; DText=$true},
    @{Row=46; E='  -0.02%  '},
    @{Row=47; D='2.846.79'; E='  -3.66%  '},
    @{Row=48; D='139.11'; E='  +0.03%  '; DText=$true},
    @{Row=49; D='0.0351'; E='  -1.96%  '; DText=$true},
    @{Row=50; D='25.86'; E='  -4.72%  '; DText=$true},
    @{Row=51; D='23.60'; E='  +5.86%  '; DText=$true}
)

foreach ($item in $changes) {
    if ($item.ContainsKey("B")) { $ws.Cells.Item($item.Row, 2).Value = $item.B }
    if ($item.ContainsKey("C")) { $ws.Cells.Item($item.Row, 3).Value = $item.C }
    if ($item.ContainsKey("D")) {
        if ($item.ContainsKey("DText")) { $ws.Cells.Item($item.Row, 4).NumberFormat = "@" }
        $ws.Cells.Item($item.Row, 4).Value = $item.D
    }
    if ($item.ContainsKey("E")) { $ws.Cells.Item($item.Row, 5).Value = $item.E }
}
